# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Somalia's case count overtook Mayotte/Kenia/Mauricio, so it jumps three
#   spots up the (sorted-by-total-cases) country list; those three rows
#   cascade down one slot each
# - Refresh several countries'/provincias' stat columns with newer counts

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Datos actualizados..." timestamp -------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 16:22"

# --- Row 50 --------------------------------------------------------------
$ws.Range("D50").Value = 113
$ws.Range("E50").Value = 4745

# --- Row 52 --------------------------------------------------------------
$ws.Range("E52").Value = 1789
$ws.Range("F52").Value = 60
$ws.Range("G52").Value = 9
$ws.Range("H52").Value = 186

# --- Row 110 ---------------------------------------------------------------
$ws.Range("B110").Value = 440
$ws.Range("C110").Value = 23
$ws.Range("E110").Value = 315

# --- Somalia jumps ahead of Mayotte/Kenia/Mauricio; rows 115-118 cascade --
# Row 115 becomes Somalia with its refreshed counts
$ws.Range("A115").Value = "Somalia"
$ws.Range("B115").Value = 390
$ws.Range("C115").Value = 62
$ws.Range("D115").Value = 8
$ws.Range("E115").Value = 366
$ws.Range("F115").Value = 2
$ws.Range("H115").Value = 16

# Row 116 becomes Mayotte (old row-115 figures)
$ws.Range("A116").Value = "Mayotte"
$ws.Range("B116").Value = 354
$ws.Range("C116").Value = 0
$ws.Range("D116").Value = 144
$ws.Range("E116").Value = 206
$ws.Range("F116").Value = 4
$ws.Range("H116").Value = 4

# Row 117 becomes Kenia (old row-116 figures)
$ws.Range("A117").Value = "Kenia"
$ws.Range("B117").Value = 343
$ws.Range("C117").Value = 7
$ws.Range("D117").Value = 98
$ws.Range("E117").Value = 231
$ws.Range("F117").Value = 2
$ws.Range("H117").Value = 14

# Row 118 becomes Mauricio (old row-117 figures)
$ws.Range("A118").Value = "Mauricio"
$ws.Range("B118").Value = 331
$ws.Range("D118").Value = 285
$ws.Range("E118").Value = 37
$ws.Range("F118").Value = 3
$ws.Range("H118").Value = 9

# --- Row 122 ---------------------------------------------------------------
$ws.Range("D122").Value = 243
$ws.Range("E122").Value = 47
$ws.Range("F122").Value = 22
